$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 779. This shifts the existing row 779
# ("2026/12/29", ...) and everything below it down by one row, turning the
# old A1:D820 range into A1:D821.
$ws.Rows.Item(779).Insert()

# Populate the freshly inserted row 779 with the new reading.
# Column A holds a literal date-like string (e.g. "2026/12/29") elsewhere in
# the sheet, not a real Excel date value, so force the cell to stay text
# (apostrophe prefix) and then strip the resulting "quote prefix" number
# format so the cell keeps the sheet's default (no explicit style index).
$ws.Range("A779").Value = "'2026/02/04"
$ws.Range("A779").ClearFormats()

$ws.Range("B779").Value = "水"
$ws.Range("C779").Value = 13
$ws.Range("D779").Value = 201
